$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (C column) date values from 45221 to 45224 (2023-10-22 -> 2023-10-25)
# for rows 2 through 5.
foreach ($row in 2..5) {
    $ws.Cells.Item($row, 3).Value = 45224
}
